# Auto-generated edit script applying the cryptos.xlsx diff
# (cryptocurrency price/volume refresh, GitHub Actions commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.398.42'
$ws.Range("E2").Value = '  +1.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.362.92'
$ws.Range("E3").Value = '  +1.04%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.89'
$ws.Range("E5").Value = '  +0.61%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.53'
$ws.Range("E6").Value = '  -1.09%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("E8").Value = '  +1.52%  '

$ws.Range("E9").Value = '  +4.97%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.69'
$ws.Range("E10").Value = '  +5.30%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.356'
$ws.Range("E12").Value = '  -0.57%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.23'
$ws.Range("E13").Value = '  +2.47%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.784.35'
$ws.Range("E14").Value = '  +1.03%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.245.25'
$ws.Range("E15").Value = '  +0.78%  '

$ws.Range("E16").Value = '  +2.69%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.352.47'
$ws.Range("E17").Value = '  +0.00%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.06'
$ws.Range("E18").Value = '  +4.26%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.34'
$ws.Range("E19").Value = '  +2.49%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '331.79'
$ws.Range("E20").Value = '  -0.91%  '

$ws.Range("E21").Value = '  +4.15%  '

$ws.Range("E22").Value = '  +0.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '63.97'

$ws.Range("E24").Value = '  +1.43%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.25%  '

$ws.Range("E26").Value = '  -2.40%  '

$ws.Range("E27").Value = '  -5.11%  '

$ws.Range("E28").Value = '  +0.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.86'
$ws.Range("E29").Value = '  +0.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0745'
$ws.Range("E30").Value = '  +1.61%  '

$ws.Range("E31").Value = '  +0.67%  '

$ws.Range("E32").Value = '  -0.02%  '

$ws.Range("E33").Value = '  -2.25%  '

$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("E36").Value = '  +0.32%  '

$ws.Range("E37").Value = '  -1.04%  '

$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '40.39'
$ws.Range("E38").Value = '  +3.26%  '

$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.427'
$ws.Range("E39").Value = '  +13.68%  '

$ws.Range("E40").Value = '  -0.86%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.69'
$ws.Range("E41").Value = '  +2.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '140.08'
$ws.Range("E42").Value = '  -4.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '287.72'
$ws.Range("E43").Value = '  +0.39%  '

$ws.Range("B44").Value = 'Polygon'
$ws.Range("C44").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.424'
$ws.Range("E44").Value = '  +10.07%  '

$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0959'
$ws.Range("E45").Value = '  +2.51%  '

$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0516'
$ws.Range("E46").Value = '  +2.43%  '

$ws.Range("E47").Value = '  +1.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.65'
$ws.Range("E48").Value = '  -2.06%  '

$ws.Range("E49").Value = '  +3.11%  '

$ws.Range("E50").Value = '  -0.50%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.69'
$ws.Range("E51").Value = '  +0.05%  '
